$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '28.011.28'
Set-TextValue $ws 'E2' '  +6.67%  '
Set-TextValue $ws 'D3' '1.741.87'
Set-TextValue $ws 'E3' '  +5.11%  '
Set-TextValue $ws 'E4' '  -0.14%  '
Set-TextValue $ws 'D5' '229.78'
Set-TextValue $ws 'E5' '  +4.80%  '
Set-TextValue $ws 'D6' '0.5455'
Set-TextValue $ws 'E6' '  +3.98%  '
Set-TextValue $ws 'D7' '1.003'
Set-TextValue $ws 'E7' '  -0.17%  '
Set-TextValue $ws 'E8' '  +4.50%  '
Set-TextValue $ws 'D9' '0.06737'
Set-TextValue $ws 'E9' '  +5.76%  '
Set-TextValue $ws 'D10' '21.90'
Set-TextValue $ws 'E10' '  +5.66%  '
Set-TextValue $ws 'D11' '0.07786'
Set-TextValue $ws 'E11' '  +0.86%  '
Set-TextValue $ws 'D12' '4.727'
Set-TextValue $ws 'E12' '  +2.76%  '
Set-TextValue $ws 'B13' 'WrappedEther'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws 'D13' '1.759.76'
Set-TextValue $ws 'E13' '  +8.12%  '
Set-TextValue $ws 'B14' 'WrappedliquidstakedEther2.0'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws 'D14' '1.981.63'
Set-TextValue $ws 'E14' '  +5.07%  '
Set-TextValue $ws 'D15' '0.6048'
Set-TextValue $ws 'E15' '  +7.04%  '
Set-TextValue $ws 'D16' '0.0₅8444'
Set-TextValue $ws 'E16' '  +2.07%  '
Set-TextValue $ws 'D17' '69.79'
Set-TextValue $ws 'E17' '  +6.50%  '
Set-TextValue $ws 'D18' '27.994.39'
Set-TextValue $ws 'E18' '  +6.61%  '
Set-TextValue $ws 'D19' '226.70'
Set-TextValue $ws 'E19' '  +17.31%  '
Set-TextValue $ws 'D20' '4.839'
Set-TextValue $ws 'E20' '  +3.11%  '
Set-TextValue $ws 'E21' '  -0.18%  '
Set-TextValue $ws 'D22' '10.98'
Set-TextValue $ws 'E22' '  +5.07%  '
Set-TextValue $ws 'D23' '6.256'
Set-TextValue $ws 'E23' '  +4.13%  '
Set-TextValue $ws 'D24' '1.003'
Set-TextValue $ws 'E24' '  -0.19%  '
Set-TextValue $ws 'D25' '147.03'
Set-TextValue $ws 'E25' '  +2.54%  '
Set-TextValue $ws 'D26' '0.1252'
Set-TextValue $ws 'E26' '  +4.08%  '
Set-TextValue $ws 'D27' '7.484'
Set-TextValue $ws 'E27' '  +2.56%  '
Set-TextValue $ws 'D28' '17.20'
Set-TextValue $ws 'E28' '  +7.88%  '
Set-TextValue $ws 'D29' '1.636'
Set-TextValue $ws 'E29' '  +8.06%  '
Set-TextValue $ws 'D30' '0.05666'
Set-TextValue $ws 'E31' '  +3.31%  '
Set-TextValue $ws 'D32' '3.715'
Set-TextValue $ws 'E32' '  +5.68%  '
Set-TextValue $ws 'D33' '3.552'
Set-TextValue $ws 'E33' '  +5.64%  '
Set-TextValue $ws 'D34' '1.664'
Set-TextValue $ws 'D35' '0.9869'
Set-TextValue $ws 'E35' '  +3.91%  '
Set-TextValue $ws 'D36' '2.867'
Set-TextValue $ws 'E36' '  +2.23%  '
Set-TextValue $ws 'D37' '2.448'
Set-TextValue $ws 'E37' '  +1.44%  '
Set-TextValue $ws 'D38' '0.5963'
Set-TextValue $ws 'E38' '  +3.23%  '
Set-TextValue $ws 'D39' '0.01686'
Set-TextValue $ws 'E39' '  +5.15%  '
Set-TextValue $ws 'D40' '6.022'
Set-TextValue $ws 'E40' '  +0.75%  '
Set-TextValue $ws 'B41' 'TrustWalletToken'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D41' '0.8476'
Set-TextValue $ws 'E41' '  +0.03%  '
Set-TextValue $ws 'B42' 'Maker'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws 'D42' '1.046.93'
Set-TextValue $ws 'E42' '  +2.74%  '
Set-TextValue $ws 'D43' '1.003'
Set-TextValue $ws 'E43' '  -0.10%  '
Set-TextValue $ws 'D44' '102.46'
Set-TextValue $ws 'E44' '  +0.47%  '
Set-TextValue $ws 'D45' '1.884.95'
Set-TextValue $ws 'E45' '  +4.90%  '
Set-TextValue $ws 'D46' '0.0₈117'
Set-TextValue $ws 'E46' '  +6.18%  '
Set-TextValue $ws 'D47' '60.28'
Set-TextValue $ws 'E47' '  +3.33%  '
Set-TextValue $ws 'D48' '8.348'
Set-TextValue $ws 'E48' '  +4.04%  '
Set-TextValue $ws 'D49' '1.013'
Set-TextValue $ws 'E49' '  +1.00%  '
Set-TextValue $ws 'D50' '0.4425'
Set-TextValue $ws 'E50' '  +1.74%  '
Set-TextValue $ws 'D51' '0.05319'
Set-TextValue $ws 'E51' '  -0.28%  '
